# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$wsNote = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsNote.Range("A1")
$oldText = $cellA1.Value()
$newText = $oldText.Replace(
    "1000 Bs = 5.14 = 20386.63 pesos",
    "1000 Bs = 5.04 = 19953.81 pesos"
)
$newText = $newText.Replace(
    "20386.63 pesos = 5.13 = 971.46 Bs",
    "19953.81 pesos = 5.0 = 939.97 Bs"
)
$cellA1.Value = $newText

# --- tasas sheet: update the tasas table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 198.318
$wsTasas.Range("O10").Value = 3957.2
$wsTasas.Range("N12").Value = 3993
$wsTasas.Range("O12").Value = 188.1
